$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.400.34"
$ws.Range("E2").Value = "  +2.20%  "
$ws.Range("D3").Value = "1.664.55"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  -0.67%  "
$ws.Range("D5").Value = "'220.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").Value = "'0.506"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "'19.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.91%  "
$ws.Range("D11").Value = "'0.0849"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "1.898.39"
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("D13").Value = "1.654.39"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "'67.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.07%  "
$ws.Range("D17").Value = "27.368.02"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").Value = "0.0₃0738"
$ws.Range("D19").Value = "'224.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "'6.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.55%  "
$ws.Range("D22").Value = "'4.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.49%  "
$ws.Range("D23").Value = "'2.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.42%  "
$ws.Range("D24").Value = "'9.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'147.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").Value = "'7.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.54%  "
$ws.Range("D28").Value = "'0.120"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").Value = "'16.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("D30").Value = "'0.0513"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D35").Value = "1.267.37"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("E39").Value = "  +1.81%  "
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").Value = "'0.813"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("D43").Value = "1.810.44"
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("E44").Value = "  -4.47%  "
$ws.Range("D45").Value = "'61.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("D46").Value = "'92.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("D47").Value = "'1.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("D49").Value = "'0.0984"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("D50").Value = "'7.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").Value = "'0.407"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.01%  "
